$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header height grew slightly: 18.75 -> 19.5
$ws.Rows(1).RowHeight = 19.5

# C2 ("RETORNO VALOR API" numeric code on the Franquia-01/SUBWAY row) was
# refreshed with a new value: 649739 -> 653839
$ws.Range("C2").Value = 653839

# B2 ("SUBWAY") was using a theme-based black font color that duplicated
# another already-present explicit-black Calibri font in the style table.
# Re-asserting the color as explicit black collapses the two equivalent
# font entries into one (matching the surrounding cells, e.g. A2).
$ws.Range("B2").Font.Color = 0
